$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.23349936394030379
$ws.Range("B1").Value = 0.23298147336571162
$ws.Range("A2").Value = -0.18764099340159124
$ws.Range("B2").Value = 0.1862926465177237
$ws.Range("A3").Value = -0.13657563870962797
$ws.Range("B3").Value = 0.13618488982467802
$ws.Range("A4").Value = -0.12818488984966514
$ws.Range("B4").Value = 0.12757270615337291
$ws.Range("A5").Value = -0.12457270616835725
$ws.Range("B5").Value = 0.12247396332411142
$ws.Range("A6").Value = -0.023179846315025898
$ws.Range("B6").Value = 0.022956069177119787
$ws.Range("A7").Value = -0.012956069213228893
$ws.Range("B7").Value = 0.012914651945274969
$ws.Range("A8").Value = -0.031957598239634954
$ws.Range("B8").Value = 0.031668317620953346
$ws.Range("A9").Value = -0.029668317639941932
$ws.Range("B9").Value = 0.029425256853057302
$ws.Range("A10").Value = -0.027425256873421233
$ws.Range("B10").Value = 0.027408278317480494
$ws.Range("A11").Value = -0.024408278340398937
$ws.Range("B11").Value = 0.024380534312975044
$ws.Range("A12").Value = -0.020880534337439638
$ws.Range("B12").Value = 0.020678498479817176
$ws.Range("A13").Value = -0.017178498505675321
$ws.Range("B13").Value = 0.017086248921816427
$ws.Range("A14").Value = -0.0090862489587566486
$ws.Range("B14").Value = 0.0090558904492672454
$ws.Range("A15").Value = -0.008055890469878868
$ws.Range("B15").Value = 0.0080362656873678162
$ws.Range("A16").Value = -0.0060362657106090012
$ws.Range("B16").Value = 0.0060038491088953982
$ws.Range("A17").Value = -0.0040038491324985159
$ws.Range("B17").Value = 0.0039999999715565338
$ws.Range("A18").Value = -0.060866368541727667
$ws.Range("B18").Value = 0.060764663039400091
$ws.Range("A19").Value = -0.056764663049903241
$ws.Range("B19").Value = 0.05603230883929422
$ws.Range("A20").Value = -0.05203230885295973
$ws.Range("B20").Value = 0.051825178499450786
$ws.Range("A21").Value = -0.0040058662946007573
$ws.Range("B21").Value = 0.0039999999854893886
$ws.Range("A22").Value = -0.045717007821753697
$ws.Range("B22").Value = 0.045501988298912721
$ws.Range("A23").Value = -0.040501988315685189
$ws.Range("B23").Value = 0.040099699145535794
$ws.Range("A24").Value = -0.020099699199198184
$ws.Range("B24").Value = 0.019999999945681246
$ws.Range("A25").Value = -0.097294117026025972
$ws.Range("B25").Value = 0.097166216167233088
$ws.Range("A26").Value = -0.094666216186096719
$ws.Range("B26").Value = 0.09450207415337708
$ws.Range("A27").Value = -0.092002074173309634
$ws.Range("B27").Value = 0.091033798632786045
$ws.Range("A28").Value = -0.089033798655783869
$ws.Range("B28").Value = 0.088370737128710353
$ws.Range("A29").Value = -0.081370737166752249
$ws.Range("B29").Value = 0.081178194826104111
$ws.Range("A30").Value = -0.021178194991379407
$ws.Range("B30").Value = 0.021025230108056903
$ws.Range("A31").Value = -0.014025230148831724
$ws.Range("B31").Value = 0.014001642118330437
$ws.Range("A32").Value = -0.004001642166405972
$ws.Range("B32").Value = 0.0039999999660942365
